$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update footer timestamp text (A1)
$ws.Range("A1").Value = "Datos actualizados a 13 de Mayo de 2020 a las 14:35"

# Row 4
$ws.Range("B4").Value = 1408823
$ws.Range("C4").Value = 187
$ws.Range("E4").Value = 1028622
$ws.Range("G4").Value = 30
$ws.Range("H4").Value = 83455

# Row 5
$ws.Range("B5").Value = 271095
$ws.Range("C5").Value = 1575
$ws.Range("D5").Value = 183227
$ws.Range("E5").Value = 60764
$ws.Range("G5").Value = 184
$ws.Range("H5").Value = 27104

# Row 11
$ws.Range("B11").Value = 173424
$ws.Range("C11").Value = 253
$ws.Range("E11").Value = 16951
$ws.Range("G11").Value = 35
$ws.Range("H11").Value = 7773

# Row 47
$ws.Range("D47").Value = 8663
$ws.Range("E47").Value = 1471
$ws.Range("F47").Value = 37
$ws.Range("G47").Value = 6
$ws.Range("H47").Value = 533

# Row 52
$ws.Range("B52").Value = 8158
$ws.Range("C52").Value = 1
$ws.Range("G52").Value = 1
$ws.Range("H52").Value = 229

# Row 77
$ws.Range("B77").Value = 2213
$ws.Range("C77").Value = 6
$ws.Range("D77").Value = 1834
$ws.Range("E77").Value = 285
$ws.Range("F77").Value = 9
$ws.Range("G77").Value = 3
$ws.Range("H77").Value = 94

# Row 193
$ws.Range("A193").Value = "Nueva Caledonia"
$ws.Range("D193").Value = 18
$ws.Range("H193").Value = 0

# Row 194
$ws.Range("A194").Value = "Belice"
$ws.Range("D194").Value = 16
$ws.Range("H194").Value = 2

# Row 198
$ws.Range("A198").Value = "Curazao"
$ws.Range("D198").Value = 14
$ws.Range("H198").Value = 1

# Row 199
$ws.Range("A199").Value = "Dominica"
$ws.Range("D199").Value = 15
$ws.Range("H199").Value = 0

# Row 215
$ws.Range("A215").Value = "San Bartolome"

# Row 216
$ws.Range("A216").Value = "Bonaire, San Eustaquio y Saba"

# Row 218
$ws.Range("A218").Value = "Lesoto"
$ws.Range("C218").Value = 1
$ws.Range("D218").Value = 0
$ws.Range("E218").Value = 1

# Row 219
$ws.Range("A219").Value = "San Pedro y Miquelon"
$ws.Range("B219").Value = 1
$ws.Range("C219").Value = 0
$ws.Range("D219").Value = 1
$ws.Range("E219").Value = 0
$ws.Range("F219").Value = 0
$ws.Range("G219").Value = 0
$ws.Range("H219").Value = 0
